$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A2:D9) were re-ordered: sorted ascending by column A (time),
# exactly as reflected in the diff.
$dataRange = $ws.Range("A2:D9")
$sortKey = $ws.Range("A2:A9")
$dataRange.Sort($sortKey, 1)
